$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Josivaldo Ferreira-Circuitos Elétricos 2"
$ws.Range("D3").Value = "João Rodrigues-CAD"
$ws.Range("E3").Value = "Andre Barros-EAP"

$ws.Range("C4").Value = "Josivaldo Ferreira-Circuitos Elétricos 2"
$ws.Range("D4").Value = "João Rodrigues-CAD"
$ws.Range("E4").Value = "Andre Barros-EAP"

$ws.Range("C6").Value = "Andre Lucca-Acionamentos"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "-"
